$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: drop the "---Principles of Computing" suffix, bump min_units to 32
$ws.Range("A7").Value = "BS in Information Systems---Technical Core---Computer Science Requirement"
$ws.Range("C7").Value = 32

# Row 8: drop the "---Fundamentals of Programming and Computer Science" suffix, bump min_units to 32
$ws.Range("A8").Value = "BS in Information Systems---Technical Core---Computer Science Requirement"
$ws.Range("C8").Value = 32

# Row 9 & 10: min_units bump to 32 (text already matches)
$ws.Range("C9").Value = 32
$ws.Range("C10").Value = 32

# New row 97: Information Security and Privacy summative course
$ws.Range("A97").Value = "BS in Information Systems---Concentration---Information Security and Privacy---Summative Course"
$ws.Range("B97").Value = "67-427"
$ws.Range("C97").Value = 9
